$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.552.66'
$ws.Range("E2").Value = '  -1.55%  '
$ws.Range("D3").Value = '3.400.79'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.51'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.34'
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D8").Value = '3.402.42'
$ws.Range("E8").Value = '  -0.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.566'
$ws.Range("E9").Value = '  -8.72%  '
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("E11").Value = '  -3.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.422'
$ws.Range("E12").Value = '  -4.46%  '
$ws.Range("D13").Value = '3.986.83'
$ws.Range("E13").Value = '  -0.54%  '
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.85'
$ws.Range("E15").Value = '  -3.98%  '
$ws.Range("E16").Value = '  -9.10%  '
$ws.Range("D17").Value = '63.645.79'
$ws.Range("E17").Value = '  -1.43%  '
$ws.Range("D18").Value = '3.404.27'
$ws.Range("E18").Value = '  +0.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.07'
$ws.Range("E19").Value = '  -4.55%  '
$ws.Range("E20").Value = '  -3.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '382.31'
$ws.Range("E21").Value = '  +1.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.74'
$ws.Range("E22").Value = '  -3.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.10'
$ws.Range("E24").Value = '  -1.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.513'
$ws.Range("E25").Value = '  -7.32%  '
$ws.Range("E26").Value = '  -3.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.68'
$ws.Range("E27").Value = '  -6.16%  '
$ws.Range("E28").Value = '  +0.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.01'
$ws.Range("E30").Value = '  -2.88%  '
$ws.Range("E31").Value = '  -7.40%  '
$ws.Range("E32").Value = '  -2.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.81'
$ws.Range("E33").Value = '  -1.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.89'
$ws.Range("E34").Value = '  -4.74%  '
$ws.Range("E35").Value = '  -7.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.39'
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.842'
$ws.Range("E37").Value = '  +9.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.81'
$ws.Range("E38").Value = '  -4.43%  '
$ws.Range("D39").Value = '2.807.67'
$ws.Range("E39").Value = '  -1.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.84'
$ws.Range("E40").Value = '  -3.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '42.92'
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0715'
$ws.Range("E42").Value = '  -6.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.36'
$ws.Range("E43").Value = '  -9.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.55'
$ws.Range("E44").Value = '  -3.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.33'
$ws.Range("E45").Value = '  -6.13%  '
$ws.Range("E46").Value = '  -3.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '326.49'
$ws.Range("E47").Value = '  +1.96%  '
$ws.Range("E48").Value = '  +8.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.03'
$ws.Range("E49").Value = '  -5.10%  '
$ws.Range("E50").Value = '  -5.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.27'
$ws.Range("E51").Value = '  -4.81%  '
